# Update column G ("K") values on Sheet1 per regenerated save_data
# (regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newK = @{
    2  = 3
    3  = 4
    4  = 3
    5  = 2
    6  = 1
    7  = 3
    8  = 3
    9  = 3
    10 = 2
    11 = 2
    12 = 1
    13 = 3
    14 = 4
    15 = 2
    16 = 4
    17 = 5
    18 = 4
    19 = 3
    20 = 1
    21 = 4
    22 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
